# Cyclically rotate the species data for rows 3, 4, 5:
#   new row 3 = old row 4
#   new row 4 = old row 5
#   new row 5 = old row 3
# Only the columns that actually vary between these rows need updating:
# A, B, D, E, F, G, H, I, Q, R, S

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "I", "Q", "R", "S")

# Capture the original values for rows 3, 4, 5 before overwriting anything.
# (Value2 gives back a plain scalar; Value can return a wrapped Variant.)
$orig = @{}
foreach ($r in 3, 4, 5) {
    $orig[$r] = @{}
    foreach ($col in $cols) {
        $orig[$r][$col] = $ws.Range("$col$r").Value2
    }
}

# new row 3 <- old row 4, new row 4 <- old row 5, new row 5 <- old row 3
$mapping = @{ 3 = 4; 4 = 5; 5 = 3 }

foreach ($destRow in 3, 4, 5) {
    $srcRow = $mapping[$destRow]
    foreach ($col in $cols) {
        $val = $orig[$srcRow][$col]

        # If the destination cell already holds exactly this value, leave it
        # alone entirely -- re-assigning would risk changing how a
        # numeric-looking value is typed (text vs number) even when nothing
        # actually needs to move.
        if ($ws.Range("$col$destRow").Value2 -eq $val) {
            continue
        }

        if ($col -eq "I") {
            # Column I ("Antal") holds small text labels ("1" or blank) even
            # though they look numeric. A bare Value2 assignment of "1" or
            # "" auto-coerces to a number / clears the cell entirely, which
            # loses the original text typing. Prefix with an apostrophe so
            # it is stored as text, the same way Excel treats a
            # quote-prefixed entry, then reset the style so the
            # quote-prefix formatting flag doesn't linger.
            if ([string]::IsNullOrEmpty($val)) {
                $ws.Range("$col$destRow").Value2 = "'"
            } else {
                $ws.Range("$col$destRow").Value2 = "'" + $val
            }
            $ws.Range("$col$destRow").Style = "Normal"
        } else {
            $ws.Range("$col$destRow").Value2 = $val
        }
    }
}
